$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.477.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.082.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.25'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.075.04'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.44'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.471'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.23'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.586.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.484.99'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.084.80'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.27'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.52'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.06'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.06'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.63%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.07'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.48'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.58'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.04'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '467.17'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.98'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +15.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0834'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0405'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.963.25'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.114'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.41'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.258'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.09%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0521'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.32'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.32%  '
